# "add Sub Category as a variable"
#
# Net effect required by the diff:
#   * Tab order becomes Bond, Alternative, Equity (was Equity, Alternative, Bond) -
#     each tab keeps its own data (Equity tab keeps equity tickers, Bond tab keeps
#     bond tickers, etc.) - only the left-to-right ordering of the tabs changes.
#   * Every sheet gains a new "Sub Category" column right after "Ticker" (i.e. a
#     new column B is inserted, pushing Name/Exchange Name/... one column right).

$wb = $excel.ActiveWorkbook

# --- 1. Reorder the tabs: Equity, Alternative, Bond -> Bond, Alternative, Equity ---
$bond = $wb.Worksheets.Item("Bond")
$bond.Move($wb.Worksheets.Item("Equity"))          # Bond becomes the first tab
$equity = $wb.Worksheets.Item("Equity")
$equity.Move($null, $wb.Worksheets.Item("Alternative"))   # Equity becomes the last tab

# --- 2. Insert the new "Sub Category" column (column B) on every sheet ---

function Add-SubCategoryColumn($sheetName, $values) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Columns.Item(2).Insert()
    $ws.Range("B1").Value = "Sub Category"
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($i + 2, 2).Value = $values[$i]
    }
}

Add-SubCategoryColumn "Bond" @("Traditional", "Traditional", "Traditional", "Traditional")
Add-SubCategoryColumn "Alternative" @("Gold", "Managed Futures", "Global Macro", "VIX Futures")
Add-SubCategoryColumn "Equity" @("Traditional", "Traditional", "Traditional", "Traditional", "Traditional", "Traditional", "Traditional")
